# Generate Report for Handoff
# A new handoff run updated the "Latest Handoff Date(time)" timestamp for
# the e633c1d3-9883-47d4-8bde-2e70b442c8e8 row (row 7) across all three
# report sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-03-21 18:39:33"

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-21 18:39:30"

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-21 18:39:33"
